# Applies the "added frameStoragePath to spec" edit:
#   - _Main: new row (UC#200) pointing at the new requirement text.
#   - 100-x-image-input: new "frameStoragePath" row documenting the
#     Minio storage path template + its char-replacement map and an
#     example, inserted above the existing "file" JSON payload block.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------
# 1. "_Main" sheet: new requirement row (row 11)
# ---------------------------------------------------------------
$main = $wb.Worksheets.Item("_Main")

# Bring formatting down from the row above (A:number style, B/C:status
# style, D:wrapped-text style) before filling in the new row's values.
$main.Range("A10:D10").Copy()
$main.Range("A11:D11").PasteSpecial($xlPasteFormats)

$main.Range("A11").Value = 200
$main.Range("B11").Value = "OPEN"
$main.Range("C11").Value = "OPEN"
$main.Range("D11").Value = "Processed image must be saved to Minio at identifier, described on 100-x-image-input"

$main.Rows(11).RowHeight = 20.4

# ---------------------------------------------------------------
# 2. "100-x-image-input" sheet: new frameStoragePath row
# ---------------------------------------------------------------
$imgInput = $wb.Worksheets.Item("100-x-image-input")

# Push the trailing "file" JSON-payload block down a row to make room.
$imgInput.Rows(14).Insert()

# The inserted row inherits row 13's formatting; columns C/D need the
# wrapped-text style used by the other long example/value cells
# (matches the style now sitting on what was row 16, now row 17).
$imgInput.Range("C17:D17").Copy()
$imgInput.Range("C14:D14").PasteSpecial($xlPasteFormats)

$imgInput.Range("B14").Value = "frameStoragePath"
$imgInput.Range("D14").Value = 'local/jpgdata/{hostname}/{source}/frame_{timestamp}_{frameNo}_{localID}.jpg
Replacement map:
":": "_#_",
"\": "/",
"//": "/",'
$imgInput.Range("C14").Value = "local/jpgdata/HARON/Camera_#_0/frame_1712372019_15_382991.jpg
local/jpgdata/HARON/C_#_/dev/github/fc/data/video.avi/frame_1712372019_15_382991.jpg"

$imgInput.Rows(14).RowHeight = 61.2

# Column C needed to widen to fit the new example path text.
$imgInput.Columns(3).ColumnWidth = 32.61

# Page got a paper size / orientation assigned along with this edit.
$imgInput.PageSetup.PaperSize = 9
$imgInput.PageSetup.Orientation = 1

# ---------------------------------------------------------------
# 3. Selection / active-tab bookkeeping, matching where the author's
#    cursor ended up: last touched _Main!D11, then moved to and left
#    off on 100-x-image-input!B14 (now the active tab).
# ---------------------------------------------------------------
$main.Activate()
$main.Range("D11").Select()

$imgInput.Activate()
$imgInput.Range("B14").Select()
